$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '68.341.48'
$ws.Cells.Item(2, 5).Value = '  +1.49%  '
$ws.Cells.Item(3, 4).Value = '3.567.17'
$ws.Cells.Item(3, 5).Value = '  +1.94%  '
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.11%  '
$ws.Cells.Item(5, 4).Value = '619.67'
$ws.Cells.Item(5, 5).Value = '  +2.61%  '
$ws.Cells.Item(6, 4).Value = '155.30'
$ws.Cells.Item(6, 5).Value = '  +3.92%  '
$ws.Cells.Item(7, 4).Value = '3.566.15'
$ws.Cells.Item(7, 5).Value = '  +1.92%  '
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 4).Value = '0.491'
$ws.Cells.Item(9, 5).Value = '  +2.29%  '
$ws.Cells.Item(10, 4).Value = '0.146'
$ws.Cells.Item(10, 5).Value = '  +5.17%  '
$ws.Cells.Item(11, 4).Value = '7.40'
$ws.Cells.Item(11, 5).Value = '  +6.58%  '
$ws.Cells.Item(12, 4).Value = '0.438'
$ws.Cells.Item(12, 5).Value = '  +3.95%  '
$ws.Cells.Item(13, 4).Value = '33.24'
$ws.Cells.Item(13, 5).Value = '  +5.68%  '
$ws.Cells.Item(14, 4).Value = '0.0000221'
$ws.Cells.Item(14, 5).Value = '  +0.91%  '
$ws.Cells.Item(15, 4).Value = '4.166.92'
$ws.Cells.Item(15, 5).Value = '  +1.84%  '
$ws.Cells.Item(16, 4).Value = '3.562.41'
$ws.Cells.Item(16, 5).Value = '  +1.87%  '
$ws.Cells.Item(17, 4).Value = '68.379.53'
$ws.Cells.Item(17, 5).Value = '  +1.67%  '
$ws.Cells.Item(18, 4).Value = '0.117'
$ws.Cells.Item(18, 5).Value = '  -0.15%  '
$ws.Cells.Item(19, 4).Value = '6.73'
$ws.Cells.Item(19, 5).Value = '  +6.15%  '
$ws.Cells.Item(20, 4).Value = '16.01'
$ws.Cells.Item(20, 5).Value = '  +6.89%  '
$ws.Cells.Item(21, 4).Value = '10.04'
$ws.Cells.Item(21, 5).Value = '  +11.93%  '
$ws.Cells.Item(22, 4).Value = '454.51'
$ws.Cells.Item(22, 5).Value = '  +1.92%  '
$ws.Cells.Item(23, 4).Value = '0.644'
$ws.Cells.Item(23, 5).Value = '  +4.13%  '
$ws.Cells.Item(24, 4).Value = '78.49'
$ws.Cells.Item(24, 5).Value = '  +1.65%  '
$ws.Cells.Item(25, 4).Value = '0.0000133'
$ws.Cells.Item(25, 5).Value = '  +2.46%  '
$ws.Cells.Item(26, 4).Value = '3.708.61'
$ws.Cells.Item(26, 5).Value = '  +1.91%  '
$ws.Cells.Item(27, 4).Value = '0.999'
$ws.Cells.Item(27, 5).Value = '  -0.13%  '
$ws.Cells.Item(28, 4).Value = '9.23'
$ws.Cells.Item(28, 5).Value = '  +11.88%  '
$ws.Cells.Item(29, 4).Value = '10.56'
$ws.Cells.Item(29, 5).Value = '  +3.77%  '
$ws.Cells.Item(30, 4).Value = '1.71'
$ws.Cells.Item(30, 5).Value = '  +10.46%  '
$ws.Cells.Item(31, 4).Value = '2.56'
$ws.Cells.Item(31, 5).Value = '  +3.46%  '
$ws.Cells.Item(32, 4).Value = '0.171'
$ws.Cells.Item(32, 5).Value = '  +4.62%  '
$ws.Cells.Item(33, 4).Value = '1.00'
$ws.Cells.Item(33, 5).Value = '  +0.09%  '
$ws.Cells.Item(34, 4).Value = '6.38'
$ws.Cells.Item(34, 5).Value = '  +4.27%  '
$ws.Cells.Item(35, 4).Value = '26.17'
$ws.Cells.Item(35, 5).Value = '  +1.95%  '
$ws.Cells.Item(36, 4).Value = '1.92'
$ws.Cells.Item(36, 5).Value = '  +4.46%  '
$ws.Cells.Item(37, 4).Value = '3.559.45'
$ws.Cells.Item(37, 5).Value = '  +2.06%  '
$ws.Cells.Item(38, 4).Value = '8.24'
$ws.Cells.Item(38, 5).Value = '  +3.57%  '
$ws.Cells.Item(39, 4).Value = '2.40'
$ws.Cells.Item(39, 5).Value = '  +9.30%  '
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  +0.03%  '
$ws.Cells.Item(41, 4).Value = '180.71'
$ws.Cells.Item(41, 5).Value = '  +3.94%  '
$ws.Cells.Item(42, 4).Value = '0.0917'
$ws.Cells.Item(42, 5).Value = '  +4.86%  '
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.10%  '
$ws.Cells.Item(44, 4).Value = '5.61'
$ws.Cells.Item(44, 5).Value = '  +3.77%  '
$ws.Cells.Item(45, 4).Value = '31.18'
$ws.Cells.Item(45, 5).Value = '  +11.95%  '
$ws.Cells.Item(46, 4).Value = '0.899'
$ws.Cells.Item(46, 5).Value = '  +2.10%  '
$ws.Cells.Item(47, 4).Value = '46.24'
$ws.Cells.Item(47, 5).Value = '  +1.88%  '
$ws.Cells.Item(48, 4).Value = '1.33'
$ws.Cells.Item(48, 5).Value = '  +3.77%  '
$ws.Cells.Item(49, 4).Value = '2.67'
$ws.Cells.Item(49, 5).Value = '  +4.57%  '
$ws.Cells.Item(50, 4).Value = '7.79'
$ws.Cells.Item(50, 5).Value = '  +3.38%  '
$ws.Cells.Item(51, 4).Value = '0.262'
$ws.Cells.Item(51, 5).Value = '  +7.63%  '
